$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.463.62"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "2.549.47"
$ws.Range("E3").Value = "  +4.69%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Formula = "'573.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.03%  "

$ws.Range("D6").Formula = "'150.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.67%  "

$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("D9").Value = "2.546.83"
$ws.Range("E9").Value = "  +4.73%  "

$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D11").Formula = "'5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").Formula = "'0.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.40%  "

$ws.Range("D14").Formula = "'28.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.28%  "

$ws.Range("D15").Value = "3.005.28"
$ws.Range("E15").Value = "  +4.58%  "

$ws.Range("D16").Value = "63.369.71"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").Value = "2.554.89"
$ws.Range("E18").Value = "  +4.73%  "

$ws.Range("D19").Formula = "'11.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.37%  "

$ws.Range("D20").Formula = "'341.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("D21").Formula = "'4.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.53%  "

$ws.Range("D22").Formula = "'6.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Formula = "'66.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "

$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("D26").Formula = "'1.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.15%  "

$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Formula = "'1.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.29%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Formula = "'8.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Formula = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Formula = "'7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.83%  "

$ws.Range("D31").Value = "0.0₃0837"
$ws.Range("E31").Value = "  +6.68%  "

$ws.Range("E32").Value = "  +3.86%  "

$ws.Range("D33").Formula = "'177.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("E34").Value = "  +9.02%  "

$ws.Range("D35").Formula = "'420.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.18%  "

$ws.Range("D36").Formula = "'0.408"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.12%  "

$ws.Range("E37").Value = "  +3.42%  "

$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +4.61%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").Formula = "'40.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("D43").Formula = "'155.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.91%  "

$ws.Range("D44").Formula = "'3.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.32%  "

$ws.Range("D45").Formula = "'21.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("E47").Value = "  +2.99%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Formula = "'0.0243"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.99%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Formula = "'0.0970"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").Formula = "'18.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.73%  "

$ws.Range("E51").Value = "  +7.18%  "
